$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "User" header column
$ws.Range("D1").Value = "User"
$ws.Range("A1:D1").Font.Bold = $true

# Existing rows 2-9 belong to user "tuan"
$ws.Range("D2:D9").Value = "tuan"

# New rows for the additional user "TT04658"
$ws.Range("A10").Value = "Tuan's World"
$ws.Range("B10").Value = "Folder"
$ws.Range("C10").Value = "C:\Users\tt04658\OneDrive - The Hartford\Documents\Github\tuans_world"
$ws.Range("D10").Value = "TT04658"

$ws.Range("A11").Value = "Downloads"
$ws.Range("B11").Value = "Folder"
$ws.Range("C11").Value = "C:\Users\TT04658\Downloads\"
$ws.Range("D11").Value = "TT04658"
$ws.Range("A11:C11").Interior.Color = $ws.Range("A6:C6").Interior.Color

$ws.Range("A12").Value = "Tuan's World"
$ws.Range("B12").Value = "VSCode"
$ws.Range("C12").Value = "C:\Users\TT04658\Documents\GitHub\vsc_tuansworld.bat"
$ws.Range("D12").Value = "TT04658"

[void]$ws.Range("C11").Select()
$ws.PageSetup.Orientation = 1
